$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Test"
$ws.Range("A4").Value = "User"
$ws.Range("C4").Value = "canbeanything"
$ws.Range("D4").Value = "test123"

$ws.Range("A4").Select()
